# Applies the cryptos-list refresh described by the commit diff.
# Re-creates each changed inline-string cell with its new value,
# forcing pure-numeric-looking strings to stay text (matching the
# original inlineStr cells) without leaving a residual cell style.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddress, $text) {
    $cell = $ws.Range($rangeAddress)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "35.471.96"
$ws.Range("E2").Value = "  +1.16%  "
$ws.Range("D3").Value = "1.904.30"
$ws.Range("E3").Value = "  +2.68%  "
$ws.Range("E4").Value = "  +0.43%  "
$ws.Range("E5").Value = "  +4.37%  "
Set-TextValue "D6" "0.633"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("E7").Value = "  +0.41%  "
Set-TextValue "D8" "41.97"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("E9").Value = "  +2.85%  "
$ws.Range("E10").Value = "  +1.18%  "
Set-TextValue "D11" "0.0996"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "2.179.59"
$ws.Range("E12").Value = "  +2.69%  "
$ws.Range("E13").Value = "  +8.25%  "
Set-TextValue "D14" "0.692"
$ws.Range("E14").Value = "  +1.52%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "1.895.23"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("B16").Value = "Polkadot"
$ws.Range("C16").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D16" "4.86"
$ws.Range("E16").Value = "  +3.62%  "
$ws.Range("D17").Value = "35.480.07"
$ws.Range("E17").Value = "  +1.25%  "
Set-TextValue "D18" "71.87"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("D19").Value = "0.0₃0821"
$ws.Range("E19").Value = "  +3.18%  "
Set-TextValue "D20" "243.22"
$ws.Range("E20").Value = "  +0.73%  "
Set-TextValue "D21" "12.56"
$ws.Range("E21").Value = "  +2.78%  "
Set-TextValue "D22" "4.86"
$ws.Range("E22").Value = "  +1.89%  "
$ws.Range("E23").Value = "  +0.37%  "
$ws.Range("E24").Value = "  +0.48%  "
Set-TextValue "D25" "172.26"
$ws.Range("E25").Value = "  +0.51%  "
Set-TextValue "D26" "2.17"
$ws.Range("E26").Value = "  +16.90%  "
Set-TextValue "D27" "8.54"
$ws.Range("E27").Value = "  +7.88%  "
Set-TextValue "D28" "17.94"
$ws.Range("E28").Value = "  +1.35%  "
$ws.Range("E29").Value = "  +0.34%  "
Set-TextValue "D30" "0.972"
$ws.Range("E30").Value = "  +23.88%  "
Set-TextValue "D31" "0.0570"
$ws.Range("E31").Value = "  +2.14%  "
$ws.Range("E32").Value = "  +2.52%  "
$ws.Range("E33").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("E35").Value = "  +7.46%  "
Set-TextValue "D36" "1.36"
$ws.Range("E36").Value = "  +9.78%  "
$ws.Range("E37").Value = "  -0.95%  "
Set-TextValue "D38" "1.11"
$ws.Range("E38").Value = "  +2.54%  "
$ws.Range("E39").Value = "  +1.50%  "
Set-TextValue "D40" "90.77"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue "D41" "0.0628"
$ws.Range("E41").Value = "  +16.38%  "
$ws.Range("B42").Value = "InjectiveProtocol"
$ws.Range("C42").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D42" "15.72"
$ws.Range("E42").Value = "  +4.84%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "1.349.67"
$ws.Range("E43").Value = "  +0.10%  "
Set-TextValue "D44" "49.55"
$ws.Range("E44").Value = "  +42.47%  "
Set-TextValue "D45" "2.36"
$ws.Range("E45").Value = "  +1.25%  "
Set-TextValue "D46" "13.02"
$ws.Range("E46").Value = "  +2.28%  "
$ws.Range("E47").Value = "  +0.65%  "
$ws.Range("E48").Value = "  -0.52%  "
$ws.Range("D50").Value = "2.088.52"
$ws.Range("E50").Value = "  +2.57%  "
Set-TextValue "D51" "0.0691"
$ws.Range("E51").Value = "  +1.43%  "
